$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price cells whose new text looks like a plain number (e.g. "245.56") must be
# forced to Text format first, otherwise Excel auto-converts the assigned string
# into a numeric value instead of keeping it as literal text (matching the source
# inlineStr cells). The style is reset back to Normal afterwards so no stray
# formatting is left behind.

$ws.Range("D2").Value = "26.495.70"
$ws.Range("E2").Value = "  +0.63%  "
$ws.Range("D3").Value = "1.727.64"
$ws.Range("E3").Value = "  +0.38%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.56"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.75%  "
$ws.Range("E6").Value = "  -0.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4808"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2668"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.48%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06221"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.36%  "
$ws.Range("D10").Value = "1.725.21"
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07153"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.40%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6154"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.521"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.83%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.16"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.22%  "
$ws.Range("E16").Value = "  -0.05%  "
$ws.Range("D17").Value = "26.500.63"
$ws.Range("E17").Value = "  +0.64%  "
$ws.Range("E18").Value = "  -0.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006932"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.29%  "
$ws.Range("E20").Value = "  +0.75%  "
$ws.Range("D21").Value = "1.946.59"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.528"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.51%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.956"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.47%  "
$ws.Range("E24").Value = "  -0.97%  "
$ws.Range("E25").Value = "  +0.92%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.36"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.72%  "
$ws.Range("E27").Value = "  +2.16%  "
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "106.88"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("E30").Value = "  -0.70%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08034"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +3.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.708"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.66%  "
$ws.Range("E33").Value = "  +2.22%  "
$ws.Range("B34").Value = "HuobiToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.615"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6360"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.78%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9947"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.11%  "
$ws.Range("B37").Value = "TrustWalletToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9315"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.05%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.096"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.46%  "
$ws.Range("B39").Value = "MXToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.426"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.14%  "
$ws.Range("B40").Value = "Quant"
$ws.Range("C40").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "105.47"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -8.33%  "
$ws.Range("B41").Value = "PaxDollar"
$ws.Range("C41").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.006"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.43%  "
$ws.Range("B42").Value = "VeChain"
$ws.Range("C42").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.01502"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.80%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.584"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.05%  "
$ws.Range("B44").Value = "TheSandbox"
$ws.Range("C44").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.3900"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.33%  "
$ws.Range("B45").Value = "Aptos"
$ws.Range("C45").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.943"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.98%  "
$ws.Range("B46").Value = "Algorand"
$ws.Range("C46").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1185"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.81%  "
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05333"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.87%  "
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.92"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.10%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.827"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.38%  "
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.270"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.51%  "
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3430"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.61%  "
